# "Add files via upload" - populate the word list (column A: 7-character
# Japanese words, column B left blank for rows 2-6; header row already
# contains "単語" / "意味").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "たいしょうせい"
$ws.Range("A3").Value = "かんしょうかい"
$ws.Range("A4").Value = "かいしょうかい"
$ws.Range("A5").Value = "かんしょうかん"
$ws.Range("A6").Value = "かいしょうかん"

# Match the saved selection state (active cell B5) from the target file.
$ws.Range("B5").Select() | Out-Null
